$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the TRUE() formula in F2:F31 with a literal text value "TRUE"
# (assigning "TRUE"/"FALSE" directly gets coerced to a boolean cell, and a
# leading apostrophe forces text but also stamps a brand-new "quote prefix"
# style onto the cell. To keep the original style untouched, stage the text
# in a scratch cell far outside the used range, then copy only the VALUE
# over with PasteSpecial so the destination keeps its existing style.)
$helper = $ws.Cells.Item(1000, 1000)

for ($r = 2; $r -le 31; $r++) {
    $helper.Value = "'TRUE"
    $helper.Copy()
    $dst = $ws.Cells.Item($r, 6)
    $dst.PasteSpecial(-4163)  # xlPasteValues
}

$helper.EntireRow.Delete()
$excel.CutCopyMode = $false

# Update the selection to reflect the new active range
$ws.Range("F2:F31").Select()
